$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: B6 becomes numeric 49.5, C6 gets a sourcing note referencing "(Rated)"
$ws.Range("B6").Value = 49.5
$ws.Range("C6").Value = 'Page 1 (loc: cropped_image:Battery Label Artwork): "Rating:+15.6V==49.5Wh (Rated)"'

# Row 9: C9 note updated to include "(Rated)"
$ws.Range("C9").Value = 'Page 1 (loc: cropped_image:Battery Label Artwork): "Rating:+15.6V==49.5Wh (Rated)"'

# Row 13: updated remarks text
$ws.Range("B13").Value = '主要語言為英文，並含部分中文、日文、德文等多語警告文字。頁碼採絕對 1-based。此為電池標籤展開圖，所有數值均直接取自標籤印刷內容。'
